$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.952.47"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "3.693.62"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'650.99"
$ws.Range("E5").Value = "  -4.23%  "
$ws.Range("D6").Value = "'161.63"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.504"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.446"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "'0.0000233"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "4.314.06"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'32.85"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "3.673.50"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "69.858.66"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "'16.16"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "'6.53"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "'10.58"
$ws.Range("E20").Value = "  +8.13%  "
$ws.Range("D21").Value = "'471.61"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'0.652"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'79.96"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "3.837.64"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'11.00"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "'9.20"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "'2.66"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("D30").Value = "'1.73"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'6.58"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'26.93"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.165"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.687.87"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").Value = "'8.51"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'2.28"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'5.92"
$ws.Range("E40").Value = "  -5.14%  "
$ws.Range("D41").Value = "'180.01"
$ws.Range("E41").Value = "  +6.97%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'0.0908"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").Value = "'0.931"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'29.58"
$ws.Range("E45").Value = "  +7.20%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'47.09"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "'2.74"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "'0.000270"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").Value = "'1.07"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "'7.86"
$ws.Range("E50").Value = "  -0.37%  "
